# Update the cached display text of the "datetimeFigureOut" date field
# from 7/9/2016 to 7/10/2016 everywhere it appears: on the slide master
# and on every slide layout's Date Placeholder shape.

$p = $ppt.ActivePresentation

function Update-DateFieldText($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "7/9/2016") {
                $tr.Text = "7/10/2016"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateFieldText $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateFieldText $layout.Shapes
}
